$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell without letting Excel
# auto-convert numeric-looking strings into real numbers, and
# without leaving the cell with a different persistent style.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '43.535.27'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '2.231.12'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.22%  '
Set-TextValue $ws.Range("D5") '270.50'
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("E6").Value = '  +11.64%  '
Set-TextValue $ws.Range("D7") '0.621'
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("E8").Value = '  +0.09%  '
Set-TextValue $ws.Range("D9") '0.622'
$ws.Range("E9").Value = '  +2.83%  '
Set-TextValue $ws.Range("D10") '46.55'
$ws.Range("E10").Value = '  +4.50%  '
Set-TextValue $ws.Range("D11") '0.0919'
$ws.Range("E11").Value = '  -1.13%  '
Set-TextValue $ws.Range("D12") '8.04'
$ws.Range("E12").Value = '  +13.66%  '
Set-TextValue $ws.Range("D13") '0.105'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '2.563.93'
$ws.Range("E14").Value = '  -0.15%  '
Set-TextValue $ws.Range("D15") '15.08'
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Value = '2.242.41'
$ws.Range("E16").Value = '  +0.41%  '
Set-TextValue $ws.Range("D17") '0.799'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '43.482.19'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -1.06%  '
Set-TextValue $ws.Range("D20") '5.99'
$ws.Range("E20").Value = '  -0.88%  '
Set-TextValue $ws.Range("D21") '70.34'
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("E22").Value = '  -2.53%  '
Set-TextValue $ws.Range("D23") '232.29'
$ws.Range("E23").Value = '  -0.12%  '
Set-TextValue $ws.Range("D24") '8.73'
$ws.Range("E24").Value = '  -5.71%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D26") '2.49'
$ws.Range("E26").Value = '  +10.53%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D27") '11.23'
$ws.Range("E27").Value = '  +3.90%  '
Set-TextValue $ws.Range("D28") '3.56'
$ws.Range("E28").Value = '  +5.69%  '
Set-TextValue $ws.Range("D29") '39.44'
$ws.Range("E29").Value = '  -3.03%  '
$ws.Range("E30").Value = '  +2.06%  '
Set-TextValue $ws.Range("D31") '172.98'
Set-TextValue $ws.Range("D32") '0.0924'
$ws.Range("E32").Value = '  +3.72%  '
$ws.Range("E33").Value = '  +0.27%  '
Set-TextValue $ws.Range("D34") '5.42'
$ws.Range("E34").Value = '  +1.10%  '
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("E36").Value = '  -4.28%  '
Set-TextValue $ws.Range("D37") '0.0349'
$ws.Range("E37").Value = '  -5.94%  '
$ws.Range("E38").Value = '  -5.67%  '
Set-TextValue $ws.Range("D39") '3.51'
$ws.Range("E39").Value = '  +15.31%  '
Set-TextValue $ws.Range("D40") '12.53'
$ws.Range("E40").Value = '  -4.57%  '
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("E42").Value = '  +6.47%  '
Set-TextValue $ws.Range("D43") '62.75'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("E44").Value = '  -3.33%  '
Set-TextValue $ws.Range("D45") '0.0986'
$ws.Range("E45").Value = '  -0.23%  '
Set-TextValue $ws.Range("D46") '8.38'
$ws.Range("E46").Value = '  -0.63%  '
Set-TextValue $ws.Range("D47") '99.77'
$ws.Range("E47").Value = '  -4.23%  '
$ws.Range("E48").Value = '  +1.38%  '
Set-TextValue $ws.Range("D49") '1.18'
$ws.Range("E49").Value = '  +1.87%  '
Set-TextValue $ws.Range("D50") '0.434'
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("E51").Value = '  -7.08%  '
